$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("A22").Value = "Max aft Xcg MAC"
$ws.Range("B22").Value = "%"
$ws.Range("C22").Value = 12.318389254421769

$ws.Range("A23").Value = "Max forward Xcg MAC"
$ws.Range("B23").Value = "%"
$ws.Range("C23").Value = 42.87685706504634
